$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12: average of column J (|S*|/n), bold
$ws.Range("J12").Formula = "=AVERAGE(J2:J11)"
$ws.Range("J12").Font.Bold = $true

# Row 14-17: summary labels and formulas
$ws.Range("A14").Value = "Average of SW(S*)/SW(OPT)"
$ws.Range("B14").Formula = "=AVERAGE(N2:N11)"

$ws.Range("A15").Value = "Average of SC(S*)/SC(OPT)"
$ws.Range("B15").Formula = "=AVERAGE(Z2:Z11)"

$ws.Range("A16").Value = "Worst of SW(S*)/SW(OPT)"
$ws.Range("B16").Formula = "=MIN(N2:N11)"

$ws.Range("A17").Value = "Worst of SC(S*)/SC(OPT)"
$ws.Range("B17").Formula = "=MAX(Z2:Z11)"

# Format B14 fully (bold, size 12, vertically centered), then propagate
# the same resolved style to B15:B17 and the row heights, to avoid
# generating throwaway/orphan style entries.
$b14 = $ws.Range("B14")
$b14.Font.Bold = $true
$b14.Font.Size = 12
$b14.VerticalAlignment = -4108
$ws.Rows(14).RowHeight = 15.6

$b14.Copy()
$ws.Range("B15:B17").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Rows("15:17").RowHeight = 15.6

# Selection shown in the sheet view after the edit
$ws.Range("A14:B17").Select() | Out-Null

# Page setup (A4, portrait) as recorded by the edit
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
